$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: generalize the 600mm-cut PVC pipe entry into a plain "PVC Pipe 1 1/2""
# item, quantity expressed as a length ("1m") instead of a piece count.
$ws.Range("B27").Value = "PVC Pipe 1 1/2"""
$ws.Range("C27").Value = "1m"

# Row 28: generalize the 180mm-cut PVC pipe entry into a plain "PVC Pipe 3/4""
# item, also measured in meters.
$ws.Range("B28").Value = "PVC Pipe 3/4"""
$ws.Range("C28").Value = "1m"

# Row 33 was a blank line in the mechanical-parts list; fill in the missing
# component noted in the commit message.
$ws.Range("B33").Value = "PVC cap 1 1/2"""
$ws.Range("C33").Value = 1

# Match the author's final cursor position.
$null = $ws.Range("E34").Select()
